$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 5200.5
$ws.Range("I2").Value = 5200.5
$ws.Range("K2").Value = 5200.5
$ws.Range("M2").Value = -5087.5
$ws.Range("H11").Value = 357.93332
$ws.Range("I11").Value = 357.93332
$ws.Range("K11").Value = 357.93332
$ws.Range("M11").Value = -217.93332
$ws.Range("H32").Value = 578.4286
$ws.Range("I32").Value = 640.5
$ws.Range("J32").Value = 553.6
$ws.Range("K32").Value = 640.5
$ws.Range("L32").Value = 553.6
$ws.Range("M32").Value = -314.5
$ws.Range("N32").Value = -1205.6
$ws.Range("H38").Value = 436.34784
$ws.Range("I38").Value = 192.19048
$ws.Range("J38").Value = 3000
$ws.Range("K38").Value = 576.5714400000001
$ws.Range("L38").Value = 9000
$ws.Range("M38").Value = -204.5714400000001
$ws.Range("N38").Value = -9744
$ws.Range("H58").Value = 78819.53999999999
$ws.Range("I58").Value = 1061
$ws.Range("J58").Value = 113378.89
$ws.Range("K58").Value = 3183
$ws.Range("L58").Value = 340136.67
$ws.Range("M58").Value = -3033
$ws.Range("N58").Value = -340436.67
$ws.Range("H69").Value = 4468.8
$ws.Range("J69").Value = 4586
$ws.Range("L69").Value = 13758
$ws.Range("N69").Value = -15506
$ws.Range("H72").Value = 4468.8
$ws.Range("J72").Value = 4586
$ws.Range("L72").Value = 41274
$ws.Range("N72").Value = -50010
$ws.Range("H111").Value = 742.2941
$ws.Range("I111").Value = 660.63635
$ws.Range("J111").Value = 892
$ws.Range("K111").Value = 1981.90905
$ws.Range("L111").Value = 2676
$ws.Range("M111").Value = 1085.09095
$ws.Range("N111").Value = -8810
$ws.Range("H125").Value = 1295.8182
$ws.Range("I125").Value = 819.7857
$ws.Range("J125").Value = 2128.875
$ws.Range("K125").Value = 7378.071300000001
$ws.Range("L125").Value = 19159.875
$ws.Range("M125").Value = -4918.071300000001
$ws.Range("N125").Value = -24079.875
$ws.Range("H129").Value = 5103199.5
$ws.Range("I129").Value = 50001580
$ws.Range("J129").Value = 1111.2046
$ws.Range("K129").Value = 150004740
$ws.Range("L129").Value = 3333.6138
$ws.Range("M129").Value = -149999740
$ws.Range("N129").Value = -13333.6138
$ws.Range("H132").Value = 6900550.5
$ws.Range("I132").Value = 8698911
$ws.Range("K132").Value = 26096733
$ws.Range("M132").Value = -26094203
$ws.Range("H133").Value = 26075
$ws.Range("J133").Value = 26075
$ws.Range("L133").Value = 26075
$ws.Range("N133").Value = -36195
$ws.Range("H137").Value = 3543.6667
$ws.Range("I137").Value = 4157.2144
$ws.Range("K137").Value = 12471.6432
$ws.Range("M137").Value = -9921.643199999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 1275.3334
$ws.Range("I4").Value = 1035.6
$ws.Range("J4").Value = 1575
$ws.Range("K4").Value = 1035.6
$ws.Range("L4").Value = 1575
$ws.Range("M4").Value = -919.5999999999999
$ws.Range("N4").Value = -1807
$ws.Range("H33").Value = 15940.625
$ws.Range("J33").Value = 19000
$ws.Range("L33").Value = 19000
$ws.Range("N33").Value = -19658
$ws.Range("H74").Value = 2270.7036
$ws.Range("I74").Value = 1600.2778
$ws.Range("J74").Value = 3611.5557
$ws.Range("K74").Value = 1600.2778
$ws.Range("L74").Value = 3611.5557
$ws.Range("M74").Value = -726.2778000000001
$ws.Range("N74").Value = -5359.5557
$ws.Range("H77").Value = 2270.7036
$ws.Range("I77").Value = 1600.2778
$ws.Range("J77").Value = 3611.5557
$ws.Range("K77").Value = 8001.389
$ws.Range("L77").Value = 18057.7785
$ws.Range("M77").Value = -3633.389
$ws.Range("N77").Value = -26793.7785
$ws.Range("H132").Value = 2081.077
$ws.Range("I132").Value = 1497.6586
$ws.Range("J132").Value = 4255.636
$ws.Range("K132").Value = 4492.9758
$ws.Range("L132").Value = 12766.908
$ws.Range("M132").Value = -1962.9758
$ws.Range("N132").Value = -17826.908

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 494.83334
$ws.Range("I2").Value = 494.83334
$ws.Range("K2").Value = 494.83334
$ws.Range("M2").Value = -381.83334
$ws.Range("H41").Value = 6217.8335
$ws.Range("I41").Value = 2871.4
$ws.Range("J41").Value = 8608.143
$ws.Range("K41").Value = 2871.4
$ws.Range("L41").Value = 8608.143
$ws.Range("M41").Value = -2443.4
$ws.Range("N41").Value = -9464.143
$ws.Range("H50").Value = 9011.200000000001
$ws.Range("J50").Value = 9011.200000000001
$ws.Range("L50").Value = 9011.200000000001
$ws.Range("N50").Value = -10261.2
$ws.Range("H51").Value = 8353.714
$ws.Range("J51").Value = 8353.714
$ws.Range("L51").Value = 8353.714
$ws.Range("N51").Value = -9825.714
$ws.Range("H60").Value = 20168.666
$ws.Range("J60").Value = 24103
$ws.Range("L60").Value = 24103
$ws.Range("N60").Value = -25125
$ws.Range("H61").Value = 8353.714
$ws.Range("J61").Value = 8353.714
$ws.Range("L61").Value = 8353.714
$ws.Range("N61").Value = -9049.714
$ws.Range("H99").Value = 3726.1765
$ws.Range("I99").Value = 2349.3333
$ws.Range("J99").Value = 5275.125
$ws.Range("K99").Value = 2349.3333
$ws.Range("L99").Value = 5275.125
$ws.Range("M99").Value = -851.3332999999998
$ws.Range("N99").Value = -8271.125
$ws.Range("H126").Value = 3726.1765
$ws.Range("I126").Value = 2349.3333
$ws.Range("J126").Value = 5275.125
$ws.Range("K126").Value = 7047.999899999999
$ws.Range("L126").Value = 15825.375
$ws.Range("M126").Value = -4577.999899999999
$ws.Range("N126").Value = -20765.375
$ws.Range("H141").Value = 24255.555
$ws.Range("J141").Value = 24255.555
$ws.Range("L141").Value = 24255.555
$ws.Range("N141").Value = -34615.555

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 283.375
$ws.Range("I14").Value = 283.375
$ws.Range("K14").Value = 850.125
$ws.Range("M14").Value = -677.125
$ws.Range("H18").Value = 812.8
$ws.Range("I18").Value = 354.83334
$ws.Range("K18").Value = 1064.50002
$ws.Range("M18").Value = -895.5000199999999
$ws.Range("H87").Value = 12116.667
$ws.Range("I87").Value = 8300
$ws.Range("K87").Value = 24900
$ws.Range("M87").Value = -23652
$ws.Range("H90").Value = 12116.667
$ws.Range("I90").Value = 8300
$ws.Range("K90").Value = 74700
$ws.Range("M90").Value = -68460
$ws.Range("H129").Value = 27424.7
$ws.Range("I129").Value = 3702.8572
$ws.Range("J129").Value = 40198
$ws.Range("K129").Value = 11108.5716
$ws.Range("L129").Value = 120594
$ws.Range("M129").Value = -6108.571599999999
$ws.Range("N129").Value = -130594

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4497.4707
$ws.Range("I70").Value = 4396.9287
$ws.Range("K70").Value = 4396.9287
$ws.Range("M70").Value = -4126.9287
$ws.Range("H73").Value = 4497.4707
$ws.Range("I73").Value = 4396.9287
$ws.Range("K73").Value = 4396.9287
$ws.Range("M73").Value = -3460.9287

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2598.9473
$ws.Range("I68").Value = 1041.4286
$ws.Range("K68").Value = 1041.4286
$ws.Range("M68").Value = -292.4286
$ws.Range("H71").Value = 2598.9473
$ws.Range("I71").Value = 1041.4286
$ws.Range("K71").Value = 5207.143
$ws.Range("M71").Value = -1463.143
$ws.Range("H132").Value = 2338.3264
$ws.Range("I132").Value = 1616.8928
$ws.Range("J132").Value = 3300.238
$ws.Range("K132").Value = 4850.678400000001
$ws.Range("L132").Value = 9900.714
$ws.Range("M132").Value = -2320.678400000001
$ws.Range("N132").Value = -14960.714

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 1714745.1
$ws.Range("I3").Value = 4000000
$ws.Range("J3").Value = 804
$ws.Range("K3").Value = 4000000
$ws.Range("L3").Value = 804
$ws.Range("M3").Value = -3999886
$ws.Range("N3").Value = -1032
